# Update worksheet data to match target report snapshot (new CIES ticket export).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('CIES-12269', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0225', 'Service Desk Perkons', 'P0225-PK-VLS-ITG', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12270', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0220', 'Service Desk Perkons', 'P0220-PK-VLS-IUN', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12271', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0207', 'Service Desk Perkons', 'P0207-PK-VLS-IBA', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12272', 'Prioridade 2', 'CÂMERAS DO P0 233 ESTÃO OFF-LINE   ', 'Service Desk Perkons', 'P0233-PK-VLS-LIN', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12273', 'Prioridade 2', ' câmeras do p0183 estão off-line', 'Service Desk Perkons', 'P0183-PK-VLS-JOA', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12274', 'Prioridade 2', 'CÂMERAS DO P0167 ESTÃO OFFLINE', 'Service Desk Perkons', 'P0167-PK-VLS-DOR', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12590', 'Prioridade 2', 'CÂMERAS DO P0 296 ESTÃO OFF-LINE ', 'Service Desk Perkons', 'P0296-PK-VLS-RBA', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12591', 'Prioridade 1', 'todas as câmeras do p0268 estão off-line  ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0268-PCL-PK-COL', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12604', 'Prioridade 1', 'todas as câmeras do p0292 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0292-PK-VLS-CAC', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12605', 'Prioridade 1', 'todas as câmeras do p0293 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0293-PK-VLS-CAC', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12606', 'Prioridade 1', 'todas as câmeras do p0163 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0163-PK-VLS-MIM', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12607', 'Prioridade 1', 'todas as câmeras do p0271 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0271-PCL-PK-ARA', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12609', 'Prioridade 1', 'CORRETIVA SITE P0215', 'Service Desk Perkons', 'P0215-PK-VLS-NOV', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12610', 'Prioridade 1', 'todas as câmeras do p0189 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0189-PK-VLS-BAI', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12611', 'Prioridade 1', 'todas as câmeras do p0192 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0192-PK-VLS-MON', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12612', 'Prioridade 1', 'todas as câmeras do p0234 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0234-PK-VLS-PED', 'NO PRAZO', 'FORA DO PRAZO'),
    @('CIES-12613', 'Prioridade 1', 'todas as câmeras do p0193 estão off-line', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0193-PK-VLS-ADN', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12616', 'Prioridade 1', 'CORRETIVA SITE P0215', 'Service Desk Perkons', 'P0215-PK-VLS-NOV', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12617', 'Prioridade 1', 'CORRETIVA SITE P0215', 'Service Desk Perkons', 'P0215-PK-VLS-NOV', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12623', 'Prioridade 1', 'CORRETIVA SITE P0249', 'Service Desk Perkons', 'P0249-PK-VLS-STE', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12627', 'Prioridade 1', 'CORRETIVA SITE P0259', 'Service Desk Perkons', 'P0259-PCL-PK-COL', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12628', 'Prioridade 1', 'CORRETIVA SITE P0291', 'Service Desk Perkons', 'P0291-PK-VLS-CAS', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12642', 'Prioridade 1', 'todas as câmeras do p0292 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0292-PK-VLS-CAC', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12643', 'Prioridade 1', 'todas as câmeras do p0169 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0169-PK-VLS-BOM', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12647', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0192', 'Service Desk Perkons', 'P0192-PK-VLS-MON', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12651', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0249', 'Service Desk Perkons', 'P0249-PK-VLS-STE', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12674', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0248', 'Service Desk Perkons', 'P0248-PK-VLS-BAR', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12675', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0176', 'Service Desk Perkons', 'P0176-PK-VLS-CBA', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12715', 'Prioridade 1', 'todas as câmeras do p0292 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0292-PK-VLS-CAC', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12716', 'Prioridade 1', 'todas as câmeras do p0099 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0099-PCL PK-MIM', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12717', 'Prioridade 1', 'todas as câmeras do p0185 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0185-PK-VLS-MIM', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12718', 'Prioridade 1', 'todas as câmeras do p0168 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0168-PK-VLS-BOM', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12719', 'Prioridade 1', 'todas as câmeras do p0205 estão off-line   ', 'CRISTIANA NEVES BATISTA MOREIRA', 'P0205-PK-VLS-BAI', 'NO PRAZO', 'NO PRAZO'),
    @('CIES-12734', 'Prioridade 1', 'TODAS AS CÂMERAS DO SITE ESTÃO OFFLINE - P0259', 'Service Desk Perkons', 'P0259-PCL-PK-COL', 'NO PRAZO', 'NO PRAZO')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
